$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4:D38").Formula = "=PROPER(CONCATENATE(C4, "" "", B4))"
$ws.Range("N4:N38").Formula = "=MID(K4,4,FIND("" "",K4)-4)"
$ws.Range("O4:O38").Formula = "=FIND("" "",K4)-4"

$ws.Range("F22").Select() | Out-Null
